# FINAL COMMIT BY SIDDHARTHA
#
# Two cell-value edits:
#   1. "UpComingBikes" sheet, C4: bike launch date now carries a concrete
#      day ("Launch Date : Mar 2024" -> "Launch Date : 30 Mar 2024").
#   2. "AccountVerification" sheet, B2: the stale Selenium/WebDriver stack
#      trace that had been pasted into the ErrorMessage cell is replaced
#      with the actual validation message the app shows
#      ("Enter a valid email or phone number").

$wb = $excel.ActiveWorkbook

$wsBikes = $wb.Worksheets.Item("UpComingBikes")
$wsBikes.Range("C4").Value = "Launch Date : 30 Mar 2024"

$wsAccount = $wb.Worksheets.Item("AccountVerification")
$wsAccount.Range("B2").Value = "Enter a valid email or phone number"
